$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.00331093771177416, 0.00249772494046121, 0.00063895289174589, 0.00540205626657889, 0.0237187058299611, 0.0015102522895812, 0.00228474064321258, 0.0278234941042074, 0.00758998586558755, 0.00234282726973493, 0.000193622088407846, 0.00360137084438593, 0.00205239413712316, 0.906674153387418, 0.000987472650880012, 0.000677677309427459, 0.000387244176815691, 0.986117296261158, 0.000251708714930199, 0.00305922899684396, 0.92200902278932, 0.00439522140685809, 0.0251515092841791),
    @(0.985962398590431, 0.990899761844831, 0.00253644935814278, 0.00261389819350592, 0.0141344124537727, 0.00141344124537727, 0.965806339187175, 0.969949851879102, 0.000619590682905106, 0.903460026719848, 0.0026332604023467, 0.000580866265223537, 0.995643503010823, 0.0786299301024261, 0.0031366778322071, 0.955641179545763, 0.000619590682905106, 0.000406606385656476, 0.996940771003156, 0.991093383933239, 0.00346583538250044, 0.96150792882452, 0.968439599589521),
    @(0.00273007144655062, 0.000193622088407846, 0.00311731562336631, 0.988034154936395, 0.00472437895715143, 0.816794779948496, 0.00460820570410672, 0, 0.00251708714930199, 0.00770615911863225, 0.000619590682905106, 0.00131663020117335, 0.000793850562472167, 0.00487927662787771, 0.862218521888977, 0.00329157550293337, 0.00532460743121575, 0.00809340329544794, 0.00197494530176002, 0.00311731562336631, 0.0691811721881232, 0.00931322245241737, 0.00302050457916239),
    @(0.00780297016283618, 0.00617654462021027, 0.993629833291382, 0.00377563072395299, 0.957073982999981, 0.178461478885511, 0.026603674947238, 0.00168451216914826, 0.989002265378434, 0.0846709392607509, 0.996514802408659, 0.994230061765446, 0.00147152787189963, 0.00939067128778051, 0.13253431951517, 0.0396925281236083, 0.993668557709063, 0.00480182779251457, 0.000832574980153736, 0.00273007144655062, 0.00447267024222123, 0.0241253122156176, 0.00311731562336631)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $values[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}
